$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.450.92"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "1.871.31"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "'1.022"
$ws.Range("E4").Value = "  +2.13%  "
$ws.Range("D5").Value = "'317.50"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'1.020"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("D7").Value = "'0.5132"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "'0.3973"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").Value = "'0.08354"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'1.111"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "'42.05"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "'6.254"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "'20.48"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.818.21"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.229"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.022"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "'0.00001108"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'91.18"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "'0.06772"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "'17.71"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "'1.020"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'5.967"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "28.516.23"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "'2.286"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").Value = "'162.26"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").Value = "2.021.97"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'20.71"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'2.365"
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("D30").Value = "'127.42"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Value = "'0.1049"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'1.036"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "'5.804"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'3.630"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "'0.02431"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'0.06494"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "'0.2189"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "'8.906"
$ws.Range("E38").Value = "  -6.27%  "
$ws.Range("D39").Value = "'1.277"
$ws.Range("E39").Value = "  +4.34%  "
$ws.Range("D40").Value = "'0.6450"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'1.182"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "'5.029"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").Value = "'11.22"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "'0.6025"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").Value = "'13.05"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("D46").Value = "'3.730"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").Value = "'1.217"
$ws.Range("E47").Value = "  -4.31%  "
$ws.Range("D48").Value = "'1.992"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'122.35"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").Value = "'1.207"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").Value = "'0.06862"
$ws.Range("E51").Value = "  -0.14%  "
